$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows (2 per year block) for "Fossil Gases" and "Fossil Liquids" categories.
# Work from the bottom of the sheet upward so row numbers used below stay valid.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(7).Insert()

# Rewrite the FuelGroup / Year / demand values for every data row (2-37); the sheet's
# used-range dimension (A1:K37) is recalculated automatically by Excel.

$ws.Range("A2").Value = "Hydrogen"
$ws.Range("B2").Value = 2030
$ws.Range("F2").Value = [double]"1.827624446796205e-05"
$ws.Range("H2").Value = [double]"4.101128240432117e-10"
$ws.Range("I2").Value = [double]"5.162994137005741e-06"

$ws.Range("A3").Value = "Methanol"
$ws.Range("B3").Value = 2030

$ws.Range("A4").Value = "Ammonia"
$ws.Range("B4").Value = 2030

$ws.Range("A5").Value = "Synthetic Gases"
$ws.Range("B5").Value = 2030

$ws.Range("A6").Value = "Biogenic Gases"
$ws.Range("B6").Value = 2030
$ws.Range("F6").Value = [double]"6.452553999631451e-06"
$ws.Range("I6").Value = [double]"8.397989397361688e-07"

$ws.Range("A7").Value = "Fossil Gases"
$ws.Range("B7").Value = 2030
$ws.Range("F7").Value = [double]"7.03991754360106e-05"
$ws.Range("I7").Value = [double]"3.691412585452613e-06"

$ws.Range("A8").Value = "Synthetic Liquids"
$ws.Range("B8").Value = 2030

$ws.Range("A9").Value = "Biogenic Liquids"
$ws.Range("B9").Value = 2030
$ws.Range("F9").Value = [double]"0.000346064562203222"
$ws.Range("H9").Value = [double]"0.0011221082647848"
$ws.Range("I9").Value = [double]"0.0001757780853667028"
$ws.Range("K9").Value = [double]"0.0005524486940619629"

$ws.Range("A10").Value = "Fossil Liquids"
$ws.Range("B10").Value = 2030
$ws.Range("F10").Value = [double]"0.003325367738937833"
$ws.Range("H10").Value = [double]"0.0102554645064325"
$ws.Range("I10").Value = [double]"0.0011018500620336"
$ws.Range("K10").Value = [double]"0.005409986315309055"

$ws.Range("A11").Value = "Biomass [Solid]"
$ws.Range("B11").Value = 2030

$ws.Range("A12").Value = "Renewable Energy Carrier"
$ws.Range("B12").Value = 2030

$ws.Range("A13").Value = "Overall Demand"
$ws.Range("B13").Value = 2030
$ws.Range("F13").Value = [double]"0.003766560275044659"
$ws.Range("H13").Value = [double]"0.01137757318133012"
$ws.Range("I13").Value = [double]"0.001287322353062497"
$ws.Range("K13").Value = [double]"0.005962435009371018"

$ws.Range("A14").Value = "Hydrogen"
$ws.Range("B14").Value = 2040
$ws.Range("F14").Value = [double]"8.727433796592085e-05"
$ws.Range("H14").Value = [double]"3.433097310043139e-08"
$ws.Range("I14").Value = [double]"7.358350436273957e-06"

$ws.Range("A15").Value = "Methanol"
$ws.Range("B15").Value = 2040

$ws.Range("A16").Value = "Ammonia"
$ws.Range("B16").Value = 2040

$ws.Range("A17").Value = "Synthetic Gases"
$ws.Range("B17").Value = 2040
$ws.Range("F17").Value = [double]"3.285836061700991e-11"
$ws.Range("I17").Value = [double]"2.22295410110442e-12"

$ws.Range("A18").Value = "Biogenic Gases"
$ws.Range("B18").Value = 2040
$ws.Range("F18").Value = [double]"7.710515658804638e-06"
$ws.Range("I18").Value = [double]"1.528587009732526e-06"

$ws.Range("A19").Value = "Fossil Gases"
$ws.Range("B19").Value = 2040
$ws.Range("F19").Value = [double]"3.887300101110755e-05"
$ws.Range("I19").Value = [double]"3.895215178152039e-06"

$ws.Range("A20").Value = "Synthetic Liquids"
$ws.Range("B20").Value = 2040

$ws.Range("A21").Value = "Biogenic Liquids"
$ws.Range("B21").Value = 2040
$ws.Range("F21").Value = [double]"0.0001594408980537171"
$ws.Range("H21").Value = [double]"0.0013660064530657"
$ws.Range("I21").Value = [double]"0.0001150800479230624"
$ws.Range("K21").Value = [double]"0.000625142156911351"

$ws.Range("A22").Value = "Fossil Liquids"
$ws.Range("B22").Value = 2040
$ws.Range("F22").Value = [double]"0.0009603169320751703"
$ws.Range("H22").Value = [double]"0.0096689867292966"
$ws.Range("I22").Value = [double]"0.000487755378613"
$ws.Range("K22").Value = [double]"0.005248161989956104"

$ws.Range("A23").Value = "Biomass [Solid]"
$ws.Range("B23").Value = 2040

$ws.Range("A24").Value = "Renewable Energy Carrier"
$ws.Range("B24").Value = 2040

$ws.Range("A25").Value = "Overall Demand"
$ws.Range("B25").Value = 2040
$ws.Range("F25").Value = [double]"0.001253615717623081"
$ws.Range("H25").Value = [double]"0.0110350275133354"
$ws.Range("I25").Value = [double]"0.0006156175813831751"
$ws.Range("K25").Value = [double]"0.005873304146867454"

$ws.Range("A26").Value = "Hydrogen"
$ws.Range("B26").Value = 2050
$ws.Range("F26").Value = [double]"0.0001213175349813206"
$ws.Range("H26").Value = [double]"5.818823410047024e-08"
$ws.Range("I26").Value = [double]"1.169551836969147e-05"

$ws.Range("A27").Value = "Methanol"
$ws.Range("B27").Value = 2050

$ws.Range("A28").Value = "Ammonia"
$ws.Range("B28").Value = 2050

$ws.Range("A29").Value = "Synthetic Gases"
$ws.Range("B29").Value = 2050
$ws.Range("F29").Value = [double]"3.847815179239901e-10"
$ws.Range("I29").Value = [double]"5.086557031515447e-11"

$ws.Range("A30").Value = "Biogenic Gases"
$ws.Range("B30").Value = 2050
$ws.Range("F30").Value = [double]"1.424394012031528e-06"
$ws.Range("I30").Value = [double]"4.451596508182899e-07"

$ws.Range("A31").Value = "Fossil Gases"
$ws.Range("B31").Value = 2050
$ws.Range("F31").Value = [double]"2.84332898097902e-06"
$ws.Range("I31").Value = [double]"1.485254891157024e-06"

$ws.Range("A32").Value = "Synthetic Liquids"
$ws.Range("B32").Value = 2050
$ws.Range("F32").Value = [double]"1.357403602652891e-12"
$ws.Range("H32").Value = [double]"5.262851685878789e-11"
$ws.Range("I32").Value = [double]"1.95160571120915e-12"
$ws.Range("K32").Value = [double]"4.283859938684607e-11"

$ws.Range("A33").Value = "Biogenic Liquids"
$ws.Range("B33").Value = 2050
$ws.Range("F33").Value = [double]"2.131509547760184e-05"
$ws.Range("H33").Value = [double]"0.0017891741436655"
$ws.Range("I33").Value = [double]"2.958439156361175e-05"
$ws.Range("K33").Value = [double]"0.0008905614270451276"

$ws.Range("A34").Value = "Fossil Liquids"
$ws.Range("B34").Value = 2050
$ws.Range("F34").Value = [double]"8.481940484330691e-05"
$ws.Range("H34").Value = [double]"0.0087685262034002"
$ws.Range("I34").Value = [double]"8.750451470732252e-05"
$ws.Range("K34").Value = [double]"0.004892386530442252"

$ws.Range("A35").Value = "Biomass [Solid]"
$ws.Range("B35").Value = 2050

$ws.Range("A36").Value = "Renewable Energy Carrier"
$ws.Range("B36").Value = 2050

$ws.Range("A37").Value = "Overall Demand"
$ws.Range("B37").Value = 2050
$ws.Range("F37").Value = [double]"0.0002317201444341614"
$ws.Range("H37").Value = [double]"0.01055775858792832"
$ws.Range("I37").Value = [double]"0.0001307148919997771"
$ws.Range("K37").Value = [double]"0.005782948000325978"
